$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.510.08"
$ws.Range("E2").Value = "  -0.69%  "
$ws.Range("D3").Value = "1.621.72"
$ws.Range("E3").Value = "  -1.45%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.42"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.77%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.523"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.97%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.05"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.40%  "
$ws.Range("E9").Value = "  +1.24%  "
$ws.Range("E10").Value = "  -0.16%  "
$ws.Range("E11").Value = "  -1.77%  "
$ws.Range("D12").Value = "1.852.65"
$ws.Range("E12").Value = "  -1.38%  "
$ws.Range("D13").Value = "1.627.99"
$ws.Range("E13").Value = "  -0.96%  "
$ws.Range("E14").Value = "  -0.29%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.548"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.54%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.13"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.54%  "
$ws.Range("D17").Value = "27.478.03"
$ws.Range("E17").Value = "  -0.78%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "229.46"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.00%  "
$ws.Range("D19").Value = "0.0₃0716"
$ws.Range("E19").Value = "  -1.28%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.52"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.12%  "
$ws.Range("E21").Value = "  +0.04%  "
$ws.Range("E22").Value = "  +3.50%  "
$ws.Range("E23").Value = "  +0.27%  "
$ws.Range("E24").Value = "  +5.60%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "148.98"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.64%  "
$ws.Range("E26").Value = "  -1.36%  "
$ws.Range("E27").Value = "  -0.91%  "
$ws.Range("E28").Value = "  +0.10%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.52"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.00%  "
$ws.Range("E30").Value = "  -0.94%  "
$ws.Range("E31").Value = "  -1.04%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.27"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.43%  "
$ws.Range("D33").Value = "1.465.90"
$ws.Range("E33").Value = "  +0.76%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.04"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.01%  "
$ws.Range("E35").Value = "  -3.02%  "
$ws.Range("E36").Value = "  -0.46%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.949"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +6.83%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.870"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.67%  "
$ws.Range("E39").Value = "  -0.55%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.551"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.51%  "
$ws.Range("E41").Value = "  +0.06%  "
$ws.Range("E42").Value = "  -2.46%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "67.15"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.52%  "
$ws.Range("E44").Value = "  -1.25%  "
$ws.Range("E45").Value = "  -1.87%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.28"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.23%  "
$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.76"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.65%  "
$ws.Range("B48").Value = "RocketPoolETH"
$ws.Range("C48").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D48").Value = "1.761.84"
$ws.Range("E48").Value = "  -1.49%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "87.16"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.50%  "
$ws.Range("E50").Value = "  -1.33%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0995"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.28%  "
